# Weekly update: insert a new daily record as row 262, pushing the
# existing historical rows (old 262..355) down to 263..356.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("262:262").Insert()

$ws.Range("A262").Value = 11
$ws.Range("B262").Value = "Vega Monumental Concepción"
$ws.Range("C262").Value = "Bíobío"
$ws.Range("D262").Value = 45205
$ws.Range("E262").Value = 8
$ws.Range("F262").Value = 100112040
$ws.Range("G262").Value = "Cilantro"
$ws.Range("H262").Value = "Sin especificar"
$ws.Range("I262").Value = "Primera"
$ws.Range("J262").Value = 100
$ws.Range("K262").Value = 5500
$ws.Range("L262").Value = 6000
$ws.Range("M262").Value = 5750
$ws.Range("N262").Value = "`$/caja 36 atados"
$ws.Range("O262").Value = "Región Metropolitana"
$ws.Range("P262").Value = 160
$ws.Range("Q262").Value = 36
$ws.Range("R262").Value = "Hortaliza"
